$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'75.176.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "'2.863.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +10.86%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'605.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.98%  "

$ws.Range("D6").Value = "'189.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.29%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.66%  "

$ws.Range("D9").Value = "'0.196"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.45%  "

$ws.Range("D10").Value = "'2.860.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.79%  "

$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("E12").Value = "  +5.14%  "

$ws.Range("E13").Value = "  +3.67%  "

$ws.Range("D14").Value = "'3.393.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.12%  "

$ws.Range("D15").Value = "'75.156.50"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  +6.28%  "

$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").Value = "'2.860.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +11.00%  "

$ws.Range("D19").Value = "'9.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.97%  "

$ws.Range("D20").Value = "'12.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.50%  "

$ws.Range("D21").Value = "'380.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.46%  "

$ws.Range("D22").Value = "'2.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("E23").Value = "  +2.56%  "

$ws.Range("D24").Value = "'6.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "

$ws.Range("D25").Value = "'71.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.70%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "'4.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("D28").Value = "'3.011.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.87%  "

$ws.Range("D29").Value = "'9.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.79%  "

$ws.Range("E30").Value = "  +12.96%  "

$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").Value = "'535.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.89%  "

$ws.Range("D33").Value = "'1.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.02%  "

$ws.Range("D34").Value = "'8.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("D35").Value = "'1.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.31%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("E37").Value = "  +2.77%  "

$ws.Range("D38").Value = "'20.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.05%  "

$ws.Range("D39").Value = "'162.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").Value = "'19.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("D41").Value = "'185.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +24.76%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "'5.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.75%  "

$ws.Range("E44").Value = "  +8.86%  "

$ws.Range("E45").Value = "  +1.87%  "

$ws.Range("D46").Value = "'1.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.41%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'40.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.81%  "

$ws.Range("D49").Value = "'0.0863"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.14%  "

$ws.Range("D50").Value = "'0.581"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.22%  "

$ws.Range("E51").Value = "  +6.23%  "
